# Applies the "Update gh-pages to output generated at 456a3b4" edit.
# Sheet "展览" (Exhibitions): two outdated events are dropped, the
# remaining events shift up, two new events are appended, vote/price
# counters are refreshed, and sheets "演出"/"本地生活"/"全部类型"
# receive matching counter refreshes.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- "展览": rewrite rows 2-41 (drop 2 old rows, shift up, add 2 new rows) ----
# row 2: 杭州·大航海时代-亚洲巡回展
$ws1.Cells.Item(2, 2).Value = '2024-09-28'
$ws1.Cells.Item(2, 3).Value = '杭州·大航海时代-亚洲巡回展'
$ws1.Cells.Item(2, 4).Value = '转塘街道转塘街道江涵路300号之江文化中心 之江文化中心'
$ws1.Cells.Item(2, 5).Value = '2024.09.28 10:00-2025.01.05 21:00'
$ws1.Cells.Item(2, 6).Value = 73
$ws1.Cells.Item(2, 7).Value = 98
$ws1.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92467'
$ws1.Cells.Item(2, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/His58jYd1726224845935.jpeg'

# row 3: 杭州·彩虹社同人ONLY——星鸟之歌
$ws1.Cells.Item(3, 2).Value = '2024-09-28'
$ws1.Cells.Item(3, 3).Value = '杭州·彩虹社同人ONLY——星鸟之歌'
$ws1.Cells.Item(3, 4).Value = '丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)'
$ws1.Cells.Item(3, 5).Value = '2024.09.28 10:00-09.28 18:00'
$ws1.Cells.Item(3, 6).Value = 118
$ws1.Cells.Item(3, 7).Value = 79
$ws1.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91543'
$ws1.Cells.Item(3, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/0zc8qiBQ1724912733257.jpeg'

# row 4: 杭州·漫展嘉年华 翼起city 三周年庆（免费漫展）
$ws1.Cells.Item(4, 2).Value = '2024-09-30'
$ws1.Cells.Item(4, 3).Value = '杭州·漫展嘉年华 翼起city 三周年庆（免费漫展）'
$ws1.Cells.Item(4, 4).Value = '南庄路与港城大道交叉口附近 空港·德信之翼'
$ws1.Cells.Item(4, 5).Value = '2024.09.30 11:00-10.02 16:00'
$ws1.Cells.Item(4, 6).Value = 617
$ws1.Cells.Item(4, 7).Value = 39
$ws1.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92245'
$ws1.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/AgFJqqVY1725953058237.png'

# row 5: 杭州·【免费漫展】五福天虹ACG CLUB国潮动漫文化节
$ws1.Cells.Item(5, 2).Value = '2024-10-01'
$ws1.Cells.Item(5, 3).Value = '杭州·【免费漫展】五福天虹ACG CLUB国潮动漫文化节'
$ws1.Cells.Item(5, 4).Value = '新塘路108号 五福天虹购物中心'
$ws1.Cells.Item(5, 5).Value = '2024.10.01 11:00-10.02 19:00'
$ws1.Cells.Item(5, 6).Value = 331
$ws1.Cells.Item(5, 7).Value = 20
$ws1.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92704'
$ws1.Cells.Item(5, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/kRV6Nv2t1726823929449.png'

# row 6: 杭州·妖漫第59届动漫游戏展（免费漫展）
$ws1.Cells.Item(6, 2).Value = '2024-10-01'
$ws1.Cells.Item(6, 3).Value = '杭州·妖漫第59届动漫游戏展（免费漫展）'
$ws1.Cells.Item(6, 4).Value = '金城路333号 加州阳光.开元广场'
$ws1.Cells.Item(6, 5).Value = '2024.10.01 11:00-10.02 20:00'
$ws1.Cells.Item(6, 6).Value = 532
$ws1.Cells.Item(6, 7).Value = 39.9
$ws1.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92656'
$ws1.Cells.Item(6, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/V3SnvMIx1726732695730.jpeg'

# row 7: 杭州·数字国风嘉年华
$ws1.Cells.Item(7, 2).Value = '2024-10-01'
$ws1.Cells.Item(7, 3).Value = '杭州·数字国风嘉年华'
$ws1.Cells.Item(7, 4).Value = '小河路与桥弄街交叉口东北50米 桥西历史文化街区'
$ws1.Cells.Item(7, 5).Value = '2024.10.01 10:00-10.03 17:00'
$ws1.Cells.Item(7, 6).Value = 1529
$ws1.Cells.Item(7, 7).Value = 75
$ws1.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92220'
$ws1.Cells.Item(7, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/8aKZ9AoH1725592557833.jpeg'

# row 8: 杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！
$ws1.Cells.Item(8, 2).Value = '2024-10-01'
$ws1.Cells.Item(8, 3).Value = '杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！'
$ws1.Cells.Item(8, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(8, 5).Value = '2024.10.01 09:30-10.03 17:00'
$ws1.Cells.Item(8, 6).Value = 10849
$ws1.Cells.Item(8, 7).Value = 75
$ws1.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90057'
$ws1.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/Zk5evnyA1722331816981.jpeg'

# row 9: 杭州·第二届次元格子动漫展嘉宾内场——吴磊
$ws1.Cells.Item(9, 2).Value = '2024-10-01'
$ws1.Cells.Item(9, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场——吴磊'
$ws1.Cells.Item(9, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(9, 5).Value = '2024.10.01 09:30-10.01 17:00'
$ws1.Cells.Item(9, 6).Value = 192
$ws1.Cells.Item(9, 7).Value = 238
$ws1.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91063'
$ws1.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/oPYHR3VE1724225214229.jpeg'

# row 10: 杭州·第二届次元格子动漫展嘉宾内场——赵乾景
$ws1.Cells.Item(10, 2).Value = '2024-10-01'
$ws1.Cells.Item(10, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场——赵乾景'
$ws1.Cells.Item(10, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(10, 5).Value = '2024.10.01 09:30-10.01 17:00'
$ws1.Cells.Item(10, 6).Value = 75
$ws1.Cells.Item(10, 7).Value = 238
$ws1.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91067'
$ws1.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/jyApvRhu1724226286635.jpeg'

# row 11: 杭州·第二届次元格子动漫展嘉宾内场—袁铭喆
$ws1.Cells.Item(11, 2).Value = '2024-10-01'
$ws1.Cells.Item(11, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场—袁铭喆'
$ws1.Cells.Item(11, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(11, 5).Value = '2024.10.01 09:30-10.01 17:00'
$ws1.Cells.Item(11, 6).Value = 121
$ws1.Cells.Item(11, 7).Value = 238
$ws1.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90716'
$ws1.Cells.Item(11, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/VC67fnAB1723620167803.jpeg'

# row 12: 杭州·第六届华盟次元动漫嘉年华
$ws1.Cells.Item(12, 2).Value = '2024-10-01'
$ws1.Cells.Item(12, 3).Value = '杭州·第六届华盟次元动漫嘉年华'
$ws1.Cells.Item(12, 4).Value = '创意路1号 中国智谷富春园区'
$ws1.Cells.Item(12, 5).Value = '2024.10.01 10:00-10.02 17:00'
$ws1.Cells.Item(12, 6).Value = 2050
$ws1.Cells.Item(12, 7).Value = 60
$ws1.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89966'
$ws1.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/d0ryUws41721962610997.jpeg'

# row 13: 杭州·萌忧 原崩铁同人only
$ws1.Cells.Item(13, 2).Value = '2024-10-01'
$ws1.Cells.Item(13, 3).Value = '杭州·萌忧 原崩铁同人only'
$ws1.Cells.Item(13, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(13, 5).Value = '2024.10.01 10:00-10.01 17:00'
$ws1.Cells.Item(13, 6).Value = 891
$ws1.Cells.Item(13, 7).Value = 55
$ws1.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90671'
$ws1.Cells.Item(13, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/36k37syr1726802875327.jpeg'

# row 14: 杭州·逆光ZERO动漫游戏展
$ws1.Cells.Item(14, 2).Value = '2024-10-01'
$ws1.Cells.Item(14, 3).Value = '杭州·逆光ZERO动漫游戏展'
$ws1.Cells.Item(14, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(14, 5).Value = '2024.10.01 10:00-10.01 17:00'
$ws1.Cells.Item(14, 6).Value = 32
$ws1.Cells.Item(14, 7).Value = 40
$ws1.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91699'
$ws1.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/koWM0hw91724885204658.jpeg'

# row 15: 建德·逆光ZERO动漫游戏展
$ws1.Cells.Item(15, 2).Value = '2024-10-02'
$ws1.Cells.Item(15, 3).Value = '建德·逆光ZERO动漫游戏展'
$ws1.Cells.Item(15, 4).Value = '南山路1号 杭州新安雷迪森酒店'
$ws1.Cells.Item(15, 5).Value = '2024.10.02 10:00-10.02 17:00'
$ws1.Cells.Item(15, 6).Value = 6
$ws1.Cells.Item(15, 7).Value = 40
$ws1.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91698'
$ws1.Cells.Item(15, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/vQITgOEI1724885714305.jpeg'

# row 16: 杭州·弹丸论破only同人展
$ws1.Cells.Item(16, 2).Value = '2024-10-02'
$ws1.Cells.Item(16, 3).Value = '杭州·弹丸论破only同人展'
$ws1.Cells.Item(16, 4).Value = '北干街道萧杭路689号 杭州时尚外滩艺术中心'
$ws1.Cells.Item(16, 5).Value = '2024.10.02 09:30-10.02 17:00'
$ws1.Cells.Item(16, 6).Value = 212
$ws1.Cells.Item(16, 7).Value = 80
$ws1.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91811'
$ws1.Cells.Item(16, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/fB9EpBgU1724852399775.png'

# row 17: 杭州·第二届次元格子动漫展嘉宾内场—紫枫儿
$ws1.Cells.Item(17, 2).Value = '2024-10-02'
$ws1.Cells.Item(17, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场—紫枫儿'
$ws1.Cells.Item(17, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(17, 5).Value = '2024.10.02 09:30-10.02 17:00'
$ws1.Cells.Item(17, 6).Value = 60
$ws1.Cells.Item(17, 7).Value = 128
$ws1.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90722'
$ws1.Cells.Item(17, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/hfph4bQx1723621998996.jpeg'

# row 18: 杭州·第二届次元格子动漫展嘉宾内场——赵成晨
$ws1.Cells.Item(18, 2).Value = '2024-10-03'
$ws1.Cells.Item(18, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场——赵成晨'
$ws1.Cells.Item(18, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(18, 5).Value = '2024.10.03 09:30-10.03 17:00'
$ws1.Cells.Item(18, 6).Value = 234
$ws1.Cells.Item(18, 7).Value = 238
$ws1.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91121'
$ws1.Cells.Item(18, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/ddmmhJEE1724313674505.jpeg'

# row 19: 杭州·创世次元第五人格同人only展
$ws1.Cells.Item(19, 2).Value = '2024-10-04'
$ws1.Cells.Item(19, 3).Value = '杭州·创世次元第五人格同人only展'
$ws1.Cells.Item(19, 4).Value = '小河路与桥弄街交叉口东北50米 桥西历史文化街区'
$ws1.Cells.Item(19, 5).Value = '2024.10.04 10:00-10.05 17:00'
$ws1.Cells.Item(19, 6).Value = 1157
$ws1.Cells.Item(19, 7).Value = 75
$ws1.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92141'
$ws1.Cells.Item(19, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/MMF3dkAw1725550270634.jpeg'

# row 20: 杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only
$ws1.Cells.Item(20, 2).Value = '2024-10-04'
$ws1.Cells.Item(20, 3).Value = '杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only'
$ws1.Cells.Item(20, 4).Value = '莫干山路987号 资辉壹方汇'
$ws1.Cells.Item(20, 5).Value = '2024.10.04 09:30-10.05 17:00'
$ws1.Cells.Item(20, 6).Value = 131
$ws1.Cells.Item(20, 7).Value = 58
$ws1.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92406'
$ws1.Cells.Item(20, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/mQh43oPd1726134932363.png'

# row 21: 杭州·2024·华彩的摔跤宴 - 木吉KAZUYA降临
$ws1.Cells.Item(21, 2).Value = '2024-10-05'
$ws1.Cells.Item(21, 3).Value = '杭州·2024·华彩的摔跤宴 - 木吉KAZUYA降临'
$ws1.Cells.Item(21, 4).Value = '莫干山路188-200号 之江饭店(莫干山路店)'
$ws1.Cells.Item(21, 5).Value = '2024.10.05 10:00-10.05 16:00'
$ws1.Cells.Item(21, 6).Value = 214
$ws1.Cells.Item(21, 7).Value = 88
$ws1.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92402'
$ws1.Cells.Item(21, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/ZylQGk1P1726033298213.png'

# row 22: 杭州·文豪野犬同人only2.0
$ws1.Cells.Item(22, 2).Value = '2024-10-05'
$ws1.Cells.Item(22, 3).Value = '杭州·文豪野犬同人only2.0'
$ws1.Cells.Item(22, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(22, 5).Value = '2024.10.05 10:00-10.05 17:00'
$ws1.Cells.Item(22, 6).Value = 697
$ws1.Cells.Item(22, 7).Value = 54
$ws1.Cells.Item(22, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92226'
$ws1.Cells.Item(22, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/nkCZvaiO1725872765608.jpeg'

# row 23: 杭州·次元幻想【玩美大舞台&吃/换谷大会】（免费活动）
$ws1.Cells.Item(23, 2).Value = '2024-10-05'
$ws1.Cells.Item(23, 3).Value = '杭州·次元幻想【玩美大舞台&吃/换谷大会】（免费活动）'
$ws1.Cells.Item(23, 4).Value = '文三路 玩美的一天沉浸式生活街区'
$ws1.Cells.Item(23, 5).Value = '2024.10.05 10:00-10.05 17:00'
$ws1.Cells.Item(23, 6).Value = 72
$ws1.Cells.Item(23, 7).Value = 30
$ws1.Cells.Item(23, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92028'
$ws1.Cells.Item(23, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/FaEB96xH1725394323651.jpeg'

# row 24: 杭州·火影同人ONLY
$ws1.Cells.Item(24, 2).Value = '2024-10-05'
$ws1.Cells.Item(24, 3).Value = '杭州·火影同人ONLY'
$ws1.Cells.Item(24, 4).Value = '金城路785号B1层 萧山人民奥莱公园'
$ws1.Cells.Item(24, 5).Value = '2024.10.05 10:00-10.05 18:00'
$ws1.Cells.Item(24, 6).Value = 226
$ws1.Cells.Item(24, 7).Value = 78
$ws1.Cells.Item(24, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92458'
$ws1.Cells.Item(24, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/5choDLVP1726713753891.png'

# row 25: 杭州·第七十五届燃梦星辰国潮国漫嘉年华-随机宅舞-让我们在燃梦相遇吧！（免费展）
$ws1.Cells.Item(25, 2).Value = '2024-10-05'
$ws1.Cells.Item(25, 3).Value = '杭州·第七十五届燃梦星辰国潮国漫嘉年华-随机宅舞-让我们在燃梦相遇吧！（免费展）'
$ws1.Cells.Item(25, 4).Value = '文一西路1888号 万达广场(余杭店)'
$ws1.Cells.Item(25, 5).Value = '2024.10.05 13:00-10.06 17:00'
$ws1.Cells.Item(25, 6).Value = 2383
$ws1.Cells.Item(25, 7).Value = 58.8
$ws1.Cells.Item(25, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92605'
$ws1.Cells.Item(25, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/6j7lylE31726278092450.jpeg'

# row 26: 杭州·第五人格同人only2.0
$ws1.Cells.Item(26, 2).Value = '2024-10-05'
$ws1.Cells.Item(26, 3).Value = '杭州·第五人格同人only2.0'
$ws1.Cells.Item(26, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(26, 5).Value = '2024.10.05 10:00-10.05 17:00'
$ws1.Cells.Item(26, 6).Value = 682
$ws1.Cells.Item(26, 7).Value = 54
$ws1.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92229'
$ws1.Cells.Item(26, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/KGPiiH2U1725873923225.jpeg'

# row 27: 杭州·首届CCPC动漫嘉年华
$ws1.Cells.Item(27, 2).Value = '2024-10-05'
$ws1.Cells.Item(27, 3).Value = '杭州·首届CCPC动漫嘉年华'
$ws1.Cells.Item(27, 4).Value = '长乐路29号五组2幢 杭州运河文化发布中心'
$ws1.Cells.Item(27, 5).Value = '2024.10.05 09:00-10.06 18:00'
$ws1.Cells.Item(27, 6).Value = 3224
$ws1.Cells.Item(27, 7).Value = 69
$ws1.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91102'
$ws1.Cells.Item(27, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/ErwKRZSH1724749999253.jpeg'

# row 28: 杭州·鸢飞鱼跃代号鸢only
$ws1.Cells.Item(28, 2).Value = '2024-10-05'
$ws1.Cells.Item(28, 3).Value = '杭州·鸢飞鱼跃代号鸢only'
$ws1.Cells.Item(28, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(28, 5).Value = '2024.10.05 09:30-10.05 17:00'
$ws1.Cells.Item(28, 6).Value = 1018
$ws1.Cells.Item(28, 7).Value = 85
$ws1.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88452'
$ws1.Cells.Item(28, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/etOXBCrl1719678030944.jpeg'

# row 29: 杭州·SK怀旧动漫展SK12
$ws1.Cells.Item(29, 2).Value = '2024-10-19'
$ws1.Cells.Item(29, 3).Value = '杭州·SK怀旧动漫展SK12'
$ws1.Cells.Item(29, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws1.Cells.Item(29, 5).Value = '2024.10.19 09:00-10.20 18:00'
$ws1.Cells.Item(29, 6).Value = 737
$ws1.Cells.Item(29, 7).Value = 65
$ws1.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90921'
$ws1.Cells.Item(29, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/PspqQq6H1723894652098.jpeg'

# row 30: 杭州·SK怀旧动漫展内场陈洁丽
$ws1.Cells.Item(30, 2).Value = '2024-10-19'
$ws1.Cells.Item(30, 3).Value = '杭州·SK怀旧动漫展内场陈洁丽'
$ws1.Cells.Item(30, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws1.Cells.Item(30, 5).Value = '2024.10.19 09:00-10.19 18:00'
$ws1.Cells.Item(30, 6).Value = 9
$ws1.Cells.Item(30, 7).Value = 239
$ws1.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91272'
$ws1.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/xGSfYzkI1724393818267.jpeg'

# row 31: 杭州·风之语 动漫游戏嘉年华
$ws1.Cells.Item(31, 2).Value = '2024-10-19'
$ws1.Cells.Item(31, 3).Value = '杭州·风之语 动漫游戏嘉年华'
$ws1.Cells.Item(31, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(31, 5).Value = '2024.10.19 10:00-10.19 17:00'
$ws1.Cells.Item(31, 6).Value = 22
$ws1.Cells.Item(31, 7).Value = 55
$ws1.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92243'
$ws1.Cells.Item(31, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/rgVATfrR1726279937106.jpeg'

# row 32: 杭州·SK怀旧动漫展内场陆二喜
$ws1.Cells.Item(32, 2).Value = '2024-10-20'
$ws1.Cells.Item(32, 3).Value = '杭州·SK怀旧动漫展内场陆二喜'
$ws1.Cells.Item(32, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws1.Cells.Item(32, 5).Value = '2024.10.20 09:00-10.20 18:00'
$ws1.Cells.Item(32, 6).Value = 18
$ws1.Cells.Item(32, 7).Value = 198
$ws1.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91274'
$ws1.Cells.Item(32, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/CmcBKga31724394153544.jpeg'

# row 33: 杭州·亿万心动国乙✘代号鸢同人only(日夜场）
$ws1.Cells.Item(33, 2).Value = '2024-10-26'
$ws1.Cells.Item(33, 3).Value = '杭州·亿万心动国乙✘代号鸢同人only(日夜场）'
$ws1.Cells.Item(33, 4).Value = '皓月路299号 诺丁山艺术中心'
$ws1.Cells.Item(33, 5).Value = '2024.10.26 10:00-10.26 21:00'
$ws1.Cells.Item(33, 6).Value = 931
$ws1.Cells.Item(33, 7).Value = 75
$ws1.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91962'
$ws1.Cells.Item(33, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/LU32zDTR1725617506119.jpeg'

# row 34: 杭州·第五人格同人ONLY
$ws1.Cells.Item(34, 2).Value = '2024-10-26'
$ws1.Cells.Item(34, 3).Value = '杭州·第五人格同人ONLY'
$ws1.Cells.Item(34, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(34, 5).Value = '2024.10.26 10:00-10.26 17:00'
$ws1.Cells.Item(34, 6).Value = 24
$ws1.Cells.Item(34, 7).Value = 60
$ws1.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92639'
$ws1.Cells.Item(34, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/MdH6yT021726714388834.jpeg'

# row 35: 杭州·第三届ICIC印象动漫节【免费活动】
$ws1.Cells.Item(35, 2).Value = '2024-11-01'
$ws1.Cells.Item(35, 3).Value = '杭州·第三届ICIC印象动漫节【免费活动】'
$ws1.Cells.Item(35, 4).Value = '五常大道1号 西溪印象城'
$ws1.Cells.Item(35, 5).Value = '2024.11.01 10:00-11.03 22:00'
$ws1.Cells.Item(35, 6).Value = 31
$ws1.Cells.Item(35, 7).Value = 20
$ws1.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92700'
$ws1.Cells.Item(35, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/AznrBWao1726813980123.jpeg'

# row 36: 杭州·BanGDream! Only同人展
$ws1.Cells.Item(36, 2).Value = '2024-11-02'
$ws1.Cells.Item(36, 3).Value = '杭州·BanGDream! Only同人展'
$ws1.Cells.Item(36, 4).Value = '石祥路与丽水北路交叉口 大运河音乐公园'
$ws1.Cells.Item(36, 5).Value = '2024.11.02 10:00-11.03 20:00'
$ws1.Cells.Item(36, 6).Value = 236
$ws1.Cells.Item(36, 7).Value = 89
$ws1.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91168'
$ws1.Cells.Item(36, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/0vTxEVyz1724222524879.jpeg'

# row 37: 杭州·星梦幻夜二次元露天派对
$ws1.Cells.Item(37, 2).Value = '2024-11-02'
$ws1.Cells.Item(37, 3).Value = '杭州·星梦幻夜二次元露天派对'
$ws1.Cells.Item(37, 4).Value = '清谷路20号 森下露营'
$ws1.Cells.Item(37, 5).Value = '2024.11.02 16:00-11.02 22:00'
$ws1.Cells.Item(37, 6).Value = 3
$ws1.Cells.Item(37, 7).Value = 158
$ws1.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92348'
$ws1.Cells.Item(37, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/fysEBwBZ1726111915726.jpeg'

# row 38: 杭州·星部落动漫嘉年华
$ws1.Cells.Item(38, 2).Value = '2024-11-02'
$ws1.Cells.Item(38, 3).Value = '杭州·星部落动漫嘉年华'
$ws1.Cells.Item(38, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(38, 5).Value = '2024.11.02 09:00-11.03 16:00'
$ws1.Cells.Item(38, 6).Value = 13
$ws1.Cells.Item(38, 7).Value = 49
$ws1.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91795'
$ws1.Cells.Item(38, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/KCwYmgHz1724908471827.jpeg'

# row 39: 杭州·第六届AP动漫游戏嘉年华
$ws1.Cells.Item(39, 2).Value = '2024-11-02'
$ws1.Cells.Item(39, 3).Value = '杭州·第六届AP动漫游戏嘉年华'
$ws1.Cells.Item(39, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(39, 5).Value = '2024.11.02 09:00-11.03 17:00'
$ws1.Cells.Item(39, 6).Value = 1328
$ws1.Cells.Item(39, 7).Value = 75
$ws1.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91994'
$ws1.Cells.Item(39, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/CttbxmHX1725519238908.jpeg'

# row 40: 杭州·HD·02动漫游戏博览会
$ws1.Cells.Item(40, 2).Value = '2024-11-09'
$ws1.Cells.Item(40, 3).Value = '杭州·HD·02动漫游戏博览会'
$ws1.Cells.Item(40, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(40, 5).Value = '2024.11.09 10:00-11.10 17:00'
$ws1.Cells.Item(40, 6).Value = 1281
$ws1.Cells.Item(40, 7).Value = 75
$ws1.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92537'
$ws1.Cells.Item(40, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/nYPayxBc1725780987986.jpeg'

# row 41: 杭州·New World动漫博览会
$ws1.Cells.Item(41, 2).Value = '2024-11-09'
$ws1.Cells.Item(41, 3).Value = '杭州·New World动漫博览会'
$ws1.Cells.Item(41, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(41, 5).Value = '2024.11.09 00:00-11.10 17:00'
$ws1.Cells.Item(41, 6).Value = 5388
$ws1.Cells.Item(41, 7).Value = 75
$ws1.Cells.Item(41, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92136'
$ws1.Cells.Item(41, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/63fEMuME1725960127951.jpeg'

# ---- "展览": refresh want-to-go counters on rows 44-49 (rows 42-43 unchanged) ----
$ws1.Cells.Item(44, 6).Value = 133
$ws1.Cells.Item(45, 6).Value = 227
$ws1.Cells.Item(46, 6).Value = 51
$ws1.Cells.Item(47, 6).Value = 12
$ws1.Cells.Item(48, 6).Value = 4075
$ws1.Cells.Item(49, 6).Value = 83

# ---- "演出": refresh want-to-go counters ----
$ws2.Cells.Item(4, 6).Value = 4101
$ws2.Cells.Item(5, 6).Value = 6
$ws2.Cells.Item(6, 6).Value = 79
$ws2.Cells.Item(12, 6).Value = 394

# ---- "本地生活": refresh want-to-go counter ----
$ws3.Cells.Item(2, 6).Value = 740

# ---- "全部类型": refresh want-to-go counters ----
$ws4.Cells.Item(2, 6).Value = 740
$ws4.Cells.Item(5, 6).Value = 73
$ws4.Cells.Item(7, 6).Value = 617
$ws4.Cells.Item(8, 6).Value = 1529
$ws4.Cells.Item(9, 6).Value = 10849
$ws4.Cells.Item(14, 6).Value = 212
$ws4.Cells.Item(15, 6).Value = 234
$ws4.Cells.Item(16, 6).Value = 1157
$ws4.Cells.Item(17, 6).Value = 131
$ws4.Cells.Item(18, 6).Value = 214
$ws4.Cells.Item(19, 6).Value = 4101
$ws4.Cells.Item(20, 6).Value = 6
$ws4.Cells.Item(21, 6).Value = 697
$ws4.Cells.Item(22, 6).Value = 72
$ws4.Cells.Item(23, 6).Value = 226
$ws4.Cells.Item(24, 6).Value = 682
$ws4.Cells.Item(25, 6).Value = 3224
$ws4.Cells.Item(26, 6).Value = 1018
$ws4.Cells.Item(27, 6).Value = 79
$ws4.Cells.Item(29, 6).Value = 737
$ws4.Cells.Item(33, 6).Value = 31
$ws4.Cells.Item(34, 6).Value = 236
$ws4.Cells.Item(35, 6).Value = 1281
$ws4.Cells.Item(38, 6).Value = 133
$ws4.Cells.Item(39, 6).Value = 227
$ws4.Cells.Item(42, 6).Value = 51
$ws4.Cells.Item(43, 6).Value = 12
$ws4.Cells.Item(44, 6).Value = 4075
$ws4.Cells.Item(49, 6).Value = 83

